{"js": "// Replace each \"three-digit \u00f7 one-digit =\" problem text with its updated value.\n// Mapping derived from the authoritative diff (old text -> new text).\nconst replacements = [\n  [\"674\u00f76=\", \"277\u00f77=\"],\n  [\"946\u00f79=\", \"968\u00f75=\"],\n  [\"915\u00f76=\", \"966\u00f74=\"],\n  [\"870\u00f75=\", \"883\u00f77=\"],\n  [\"179\u00f74=\", \"346\u00f79=\"],\n  [\"637\u00f74=\", \"381\u00f76=\"],\n  [\"859\u00f77=\", \"183\u00f75=\"],\n  [\"489\u00f78=\", \"166\u00f78=\"],\n  [\"314\u00f75=\", \"298\u00f73=\"],\n  [\"163\u00f76=\", \"885\u00f78=\"],\n  [\"628\u00f72=\", \"883\u00f79=\"],\n  [\"303\u00f72=\", \"522\u00f75=\"],\n  [\"749\u00f77=\", \"567\u00f76=\"],\n  [\"986\u00f78=\", \"216\u00f76=\"],\n  [\"857\u00f75=\", \"951\u00f77=\"],\n  [\"232\u00f79=\", \"514\u00f78=\"],\n  [\"862\u00f74=\", \"216\u00f77=\"],\n  [\"435\u00f74=\", \"360\u00f79=\"],\n  [\"886\u00f77=\", \"639\u00f74=\"],\n  [\"903\u00f75=\", \"367\u00f75=\"],\n  [\"397\u00f79=\", \"277\u00f72=\"],\n  [\"763\u00f74=\", \"608\u00f75=\"],\n  [\"304\u00f74=\", \"631\u00f74=\"],\n  [\"140\u00f76=\", \"133\u00f72=\"],\n  [\"631\u00f72=\", \"244\u00f79=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"three-digit \u00f7 one-digit =\" problem text with its updated value.\n# Mapping derived from the authoritative diff (old text -> new text).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"674\u00f76=\"; New = \"277\u00f77=\" }\n    @{ Old = \"946\u00f79=\"; New = \"968\u00f75=\" }\n    @{ Old = \"915\u00f76=\"; New = \"966\u00f74=\" }\n    @{ Old = \"870\u00f75=\"; New = \"883\u00f77=\" }\n    @{ Old = \"179\u00f74=\"; New = \"346\u00f79=\" }\n    @{ Old = \"637\u00f74=\"; New = \"381\u00f76=\" }\n    @{ Old = \"859\u00f77=\"; New = \"183\u00f75=\" }\n    @{ Old = \"489\u00f78=\"; New = \"166\u00f78=\" }\n    @{ Old = \"314\u00f75=\"; New = \"298\u00f73=\" }\n    @{ Old = \"163\u00f76=\"; New = \"885\u00f78=\" }\n    @{ Old = \"628\u00f72=\"; New = \"883\u00f79=\" }\n    @{ Old = \"303\u00f72=\"; New = \"522\u00f75=\" }\n    @{ Old = \"749\u00f77=\"; New = \"567\u00f76=\" }\n    @{ Old = \"986\u00f78=\"; New = \"216\u00f76=\" }\n    @{ Old = \"857\u00f75=\"; New = \"951\u00f77=\" }\n    @{ Old = \"232\u00f79=\"; New = \"514\u00f78=\" }\n    @{ Old = \"862\u00f74=\"; New = \"216\u00f77=\" }\n    @{ Old = \"435\u00f74=\"; New = \"360\u00f79=\" }\n    @{ Old = \"886\u00f77=\"; New = \"639\u00f74=\" }\n    @{ Old = \"903\u00f75=\"; New = \"367\u00f75=\" }\n    @{ Old = \"397\u00f79=\"; New = \"277\u00f72=\" }\n    @{ Old = \"763\u00f74=\"; New = \"608\u00f75=\" }\n    @{ Old = \"304\u00f74=\"; New = \"631\u00f74=\" }\n    @{ Old = \"140\u00f76=\"; New = \"133\u00f72=\" }\n    @{ Old = \"631\u00f72=\"; New = \"244\u00f79=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $r.Old,   # FindText\n        $true,    # MatchCase\n        $false,   # MatchWholeWord\n        $false,   # MatchWildcards\n        $false,   # MatchSoundsLike\n        $false,   # MatchAllWordForms\n        $true,    # Forward\n        1,        # Wrap (wdFindContinue)\n        $false,   # Format\n        $r.New,   # ReplaceWith\n        2         # Replace (wdReplaceAll)\n    )\n}\n"}
